# Weekly update: add 4 new Cereza price observations (2022-01-21) to the
# "Feria Lagunitas de Puerto Montt" sheet, inserted right before the
# existing row 48 (pushing the existing rows 48-68 down to 52-72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at row 48 (shifts old rows 48:68 down to 52:72).
$ws.Range("A48:A51").EntireRow.Insert()

# Fixed (unchanging) column values shared by every Cereza row on this sheet.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId  = 100103
$producto    = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria   = "Cereza"

# New row data: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax,
# PrecioProm, Unidad, Origen, Precio$/Kg, Kg/unidad
$newRows = @(
    @(44582, "Lapins",  "Primera", 500, 10000, 11000, 10500, "`$/bandeja 10 kilos", "Provincia de Curicó", 1050, 10),
    @(44582, "Lapins",  "Segunda", 250,  7500,  7500,  7500, "`$/bandeja 10 kilos", "Provincia de Curicó",  750, 10),
    @(44582, "Santina", "Primera", 400, 10000, 11000, 10500, "`$/bandeja 10 kilos", "Provincia de Curicó", 1050, 10),
    @(44582, "Santina", "Segunda", 200,  7500,  7500,  7500, "`$/bandeja 10 kilos", "Provincia de Curicó",  750, 10)
)

$r = 48
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2  = $mercadoId
    $ws.Cells.Item($r, 2).Value2  = $mercado
    $ws.Cells.Item($r, 3).Value2  = $region
    $ws.Cells.Item($r, 4).Value2  = $row[0]
    $ws.Cells.Item($r, 5).Value2  = $codreg
    $ws.Cells.Item($r, 6).Value2  = $tipo
    $ws.Cells.Item($r, 7).Value2  = $productoId
    $ws.Cells.Item($r, 8).Value2  = $producto
    $ws.Cells.Item($r, 9).Value2  = $categoriaId
    $ws.Cells.Item($r, 10).Value2 = $categoria
    $ws.Cells.Item($r, 11).Value2 = $row[1]
    $ws.Cells.Item($r, 12).Value2 = $row[2]
    $ws.Cells.Item($r, 13).Value2 = $row[3]
    $ws.Cells.Item($r, 14).Value2 = $row[4]
    $ws.Cells.Item($r, 15).Value2 = $row[5]
    $ws.Cells.Item($r, 16).Value2 = $row[6]
    $ws.Cells.Item($r, 17).Value2 = $row[7]
    $ws.Cells.Item($r, 18).Value2 = $row[8]
    $ws.Cells.Item($r, 19).Value2 = $row[9]
    $ws.Cells.Item($r, 20).Value2 = $row[10]
    $r = $r + 1
}

Write-Output "Done. New dimension rows inserted; sheet now has $($ws.UsedRange.Rows.Count) used rows."
